$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6
$ws.Cells.Item(6, 7).Value = 2.2
$ws.Cells.Item(6, 9).Value = 3.05
$ws.Cells.Item(6, 15).Value = 1.88
$ws.Cells.Item(6, 16).Value = 1.39
$ws.Cells.Item(6, 17).Value = 2.85
$ws.Cells.Item(6, 20).Value = 7.8
$ws.Cells.Item(6, 22).Value = 9.5
$ws.Cells.Item(6, 23).Value = 23
$ws.Cells.Item(6, 24).Value = 18.5
$ws.Cells.Item(6, 27).Value = 7
$ws.Cells.Item(6, 32).Value = 11.5
$ws.Cells.Item(6, 33).Value = 40
$ws.Cells.Item(6, 34).Value = 27
$ws.Cells.Item(6, 35).Value = 35

# Row 10
$ws.Cells.Item(10, 10).Value = 1.08
$ws.Cells.Item(10, 11).Value = 8

# Row 13
$ws.Cells.Item(13, 7).Value = 1.55
$ws.Cells.Item(13, 10).Value = 1.04
$ws.Cells.Item(13, 11).Value = 13
$ws.Cells.Item(13, 14).Value = 1.73
$ws.Cells.Item(13, 15).Value = 2.08
$ws.Cells.Item(13, 30).Value = 15

# Row 14
$ws.Cells.Item(14, 7).Value = 3.3
$ws.Cells.Item(14, 8).Value = 3.25
$ws.Cells.Item(14, 9).Value = 2.15
$ws.Cells.Item(14, 10).Value = 1.08
$ws.Cells.Item(14, 11).Value = 8
$ws.Cells.Item(14, 22).Value = 12
$ws.Cells.Item(14, 23).Value = 34
$ws.Cells.Item(14, 30).Value = 7

# Row 18
$ws.Cells.Item(18, 9).Value = 3.8
$ws.Cells.Item(18, 11).Value = 10
$ws.Cells.Item(18, 12).Value = 1.3
$ws.Cells.Item(18, 13).Value = 3.4
$ws.Cells.Item(18, 14).Value = 2
$ws.Cells.Item(18, 15).Value = 1.8
$ws.Cells.Item(18, 22).Value = 8.5
$ws.Cells.Item(18, 26).Value = 10

# Row 19
$ws.Cells.Item(19, 12).Value = 1.29
$ws.Cells.Item(19, 13).Value = 3.5
$ws.Cells.Item(19, 14).Value = 1.93
$ws.Cells.Item(19, 15).Value = 1.88

# Row 21
$ws.Cells.Item(21, 7).Value = 2.02
$ws.Cells.Item(21, 8).Value = 3.2
$ws.Cells.Item(21, 9).Value = 3.65
$ws.Cells.Item(21, 12).Value = 1.38
$ws.Cells.Item(21, 13).Value = 2.82
$ws.Cells.Item(21, 14).Value = 2.1
$ws.Cells.Item(21, 15).Value = 1.65
$ws.Cells.Item(21, 16).Value = 1.5
$ws.Cells.Item(21, 17).Value = 2.42
$ws.Cells.Item(21, 18).Value = 1.88
$ws.Cells.Item(21, 19).Value = 1.82
$ws.Cells.Item(21, 20).Value = 6.6
$ws.Cells.Item(21, 21).Value = 9
$ws.Cells.Item(21, 22).Value = 8.5
$ws.Cells.Item(21, 23).Value = 18
$ws.Cells.Item(21, 24).Value = 17
$ws.Cells.Item(21, 27).Value = 6.2
$ws.Cells.Item(21, 28).Value = 15.5
$ws.Cells.Item(21, 29).Value = 80
$ws.Cells.Item(21, 30).Value = 9.25
$ws.Cells.Item(21, 31).Value = 18.5
$ws.Cells.Item(21, 32).Value = 12.5
$ws.Cells.Item(21, 33).Value = 55
$ws.Cells.Item(21, 34).Value = 37
$ws.Cells.Item(21, 35).Value = 45

# Row 22
$ws.Cells.Item(22, 7).Value = 1.78
$ws.Cells.Item(22, 8).Value = 3.45
$ws.Cells.Item(22, 9).Value = 4.4
$ws.Cells.Item(22, 10).Value = 1.05
$ws.Cells.Item(22, 11).Value = 7.6
$ws.Cells.Item(22, 12).Value = 1.26
$ws.Cells.Item(22, 13).Value = 3.5
$ws.Cells.Item(22, 14).Value = 1.78
$ws.Cells.Item(22, 15).Value = 1.93
$ws.Cells.Item(22, 17).Value = 2.75
$ws.Cells.Item(22, 18).Value = 1.7
$ws.Cells.Item(22, 19).Value = 2.05
$ws.Cells.Item(22, 20).Value = 7.9
$ws.Cells.Item(22, 21).Value = 9.25
$ws.Cells.Item(22, 22).Value = 7.9
$ws.Cells.Item(22, 23).Value = 15.5
$ws.Cells.Item(22, 24).Value = 13
$ws.Cells.Item(22, 25).Value = 22
$ws.Cells.Item(22, 26).Value = 7.6
$ws.Cells.Item(22, 27).Value = 6.7
$ws.Cells.Item(22, 28).Value = 13.5
$ws.Cells.Item(22, 29).Value = 55
$ws.Cells.Item(22, 31).Value = 26
$ws.Cells.Item(22, 32).Value = 14
$ws.Cells.Item(22, 33).Value = 75
$ws.Cells.Item(22, 34).Value = 40
$ws.Cells.Item(22, 35).Value = 40
$ws.Cells.Item(22, 36).Value = 400

# Row 27
$ws.Cells.Item(27, 7).Value = 2.6
$ws.Cells.Item(27, 8).Value = 2.9
$ws.Cells.Item(27, 9).Value = 2.8
$ws.Cells.Item(27, 10).Value = 1.13
$ws.Cells.Item(27, 11).Value = 4.45
$ws.Cells.Item(27, 12).Value = 1.6
$ws.Cells.Item(27, 13).Value = 2.05
$ws.Cells.Item(27, 14).Value = 2.72
$ws.Cells.Item(27, 15).Value = 1.35
$ws.Cells.Item(27, 20).Value = 5.7
$ws.Cells.Item(27, 21).Value = 10.75
$ws.Cells.Item(27, 22).Value = 11.5
$ws.Cells.Item(27, 23).Value = 29
$ws.Cells.Item(27, 24).Value = 32
$ws.Cells.Item(27, 25).Value = 60
$ws.Cells.Item(27, 26).Value = 4.75
$ws.Cells.Item(27, 27).Value = 6
$ws.Cells.Item(27, 30).Value = 5.9
$ws.Cells.Item(27, 31).Value = 11.75
$ws.Cells.Item(27, 32).Value = 12
$ws.Cells.Item(27, 33).Value = 35
$ws.Cells.Item(27, 34).Value = 35
$ws.Cells.Item(27, 35).Value = 65

# Row 28
$ws.Cells.Item(28, 7).Value = 2.18
$ws.Cells.Item(28, 9).Value = 3.15
$ws.Cells.Item(28, 13).Value = 2.42
$ws.Cells.Item(28, 14).Value = 2.25
$ws.Cells.Item(28, 18).Value = 2
$ws.Cells.Item(28, 19).Value = 1.65
$ws.Cells.Item(28, 20).Value = 6
$ws.Cells.Item(28, 21).Value = 9.25
$ws.Cells.Item(28, 22).Value = 9.5
$ws.Cells.Item(28, 23).Value = 20
$ws.Cells.Item(28, 24).Value = 21
$ws.Cells.Item(28, 26).Value = 7.3
$ws.Cells.Item(28, 28).Value = 18.5
$ws.Cells.Item(28, 29).Value = 110
$ws.Cells.Item(28, 30).Value = 7.7
$ws.Cells.Item(28, 31).Value = 15
$ws.Cells.Item(28, 32).Value = 12
$ws.Cells.Item(28, 33).Value = 40
$ws.Cells.Item(28, 34).Value = 35
$ws.Cells.Item(28, 35).Value = 50

# Row 29
$ws.Cells.Item(29, 10).Value = 1.03
$ws.Cells.Item(29, 11).Value = 17
